$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 276 - this pushes the existing rows 276:294
# down to 277:295 (matching the dimension growing from R294 to R295).
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new weekly record.
$ws.Range("A276").Value = 9
$ws.Range("B276").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C276").Value = "Metropolitana"
$ws.Range("D276").Value = 44746
$ws.Range("E276").Value = 13
$ws.Range("F276").Value = 100112043
$ws.Range("G276").Value = "Pepino ensalada"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 70
$ws.Range("K276").Value = 16000
$ws.Range("L276").Value = 18000
$ws.Range("M276").Value = 17000
$ws.Range("N276").Value = "$/caja 60 unidades"
$ws.Range("O276").Value = "Región de Arica y Parinacota"
$ws.Range("P276").Value = 283
$ws.Range("Q276").Value = 60
$ws.Range("R276").Value = "Hortaliza"
